$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# Row 1 (headers): duplicate A1:J1 into L1:U1 (column K stays empty)
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "Algoritmo"
$ws.Range("M1").Value = "y"
$ws.Range("N1").Value = "n_estimators"
$ws.Range("O1").Value = "learning_rate"
$ws.Range("P1").Value = "max_depth"
$ws.Range("Q1").Value = "min_samples_split"
$ws.Range("R1").Value = "min_samples_leaf"
$ws.Range("S1").Value = "rmse"
$ws.Range("T1").Value = "mape"
$ws.Range("U1").Value = "tiempo (min)"

# ---------------------------------------------------------------------------
# Row 2: new "Precio_m2" model-tuning result block in L2:U2
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = "GradientBoostingRegressor"
$ws.Range("M2").Value = "Precio_m2"
$ws.Range("N2").Value = 20
$ws.Range("O2").Value = 0.1
$ws.Range("P2").Value = 20
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 15
$ws.Range("S2").NumberFormat = "#,##0.00"
$ws.Range("S2").Value = 91837620.141484603
$ws.Range("T2").Value = 0.17870844335130001
$ws.Range("U2").Value = 9.5

# ---------------------------------------------------------------------------
# Row 3: new "Precio" model-tuning result (A3:J3)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "GradientBoostingRegressor"
$ws.Range("B3").Value = "Precio"
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 0.1
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 15
$ws.Range("H3").NumberFormat = "#,##0.00"
$ws.Range("H3").Value = 92202532.692044497
$ws.Range("I3").Value = 0.18619268867327099
$ws.Range("J3").Value = 8
$ws.Range("S3").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------------
# Row 4: new "Precio" model-tuning result (A4:J4)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "GradientBoostingRegressor"
$ws.Range("B4").Value = "Precio"
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 0.1
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 15
$ws.Range("H4").NumberFormat = "#,##0.00"
$ws.Range("H4").Value = 93065842.771723807
$ws.Range("I4").Value = 0.189967700342627
$ws.Range("J4").Value = 8

# ---------------------------------------------------------------------------
# Row 5: new "Precio" model-tuning result (A5:J5)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "GradientBoostingRegressor"
$ws.Range("B5").Value = "Precio"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 0.1
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 15
$ws.Range("H5").NumberFormat = "#,##0.00"
$ws.Range("H5").Value = 85558007.302483395
$ws.Range("I5").Value = 0.155435134275934
$ws.Range("J5").Value = 23
$ws.Range("S5").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------------
# Column widths for the new columns (L mirrors column A's "Algoritmo" width,
# S mirrors column H's "rmse" width). The runtime quantizes ColumnWidth to
# 1/6-character steps, so these inputs are chosen as the closest achievable
# approximation of the real Excel bestFit widths (25.28515625 / 12.7109375).
# ---------------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 24.5
$ws.Columns.Item(19).ColumnWidth = 11.9

# ---------------------------------------------------------------------------
# Final selection, matching the saved state in the workbook
# ---------------------------------------------------------------------------
$sel = $ws.Range("D5").Select()
